$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "335.40"
Set-TextValue $ws.Range("E2") "1.98%"
Set-TextValue $ws.Range("D3") "43.83"
Set-TextValue $ws.Range("E3") "6.49%"
Set-TextValue $ws.Range("D4") "5.771"
Set-TextValue $ws.Range("E4") "2.76%"
Set-TextValue $ws.Range("D5") "0.08327"
Set-TextValue $ws.Range("E5") "1.38%"
Set-TextValue $ws.Range("D6") "8.848"
Set-TextValue $ws.Range("E6") "1.14%"
Set-TextValue $ws.Range("D7") "4.520"
Set-TextValue $ws.Range("E7") "0.63%"
Set-TextValue $ws.Range("D8") "1.969"
Set-TextValue $ws.Range("E8") "-1.91%"
Set-TextValue $ws.Range("E9") "-1.90%"
Set-TextValue $ws.Range("D10") "0.9430"
Set-TextValue $ws.Range("E10") "2.44%"
Set-TextValue $ws.Range("D11") "0.1245"
Set-TextValue $ws.Range("E11") "-2.92%"
Set-TextValue $ws.Range("D12") "0.1960"
Set-TextValue $ws.Range("E12") "0.43%"
Set-TextValue $ws.Range("D13") "0.09958"
Set-TextValue $ws.Range("E13") "6.30%"
Set-TextValue $ws.Range("D14") "0.04570"
Set-TextValue $ws.Range("E14") "17.31%"
Set-TextValue $ws.Range("D15") "0.1067"
Set-TextValue $ws.Range("E15") "0.65%"
Set-TextValue $ws.Range("D16") "0.001298"
Set-TextValue $ws.Range("E16") "-0.81%"
Set-TextValue $ws.Range("D17") "0.005974"
Set-TextValue $ws.Range("E17") "-2.09%"
Set-TextValue $ws.Range("D18") "3.500"
Set-TextValue $ws.Range("E18") "1.58%"
Set-TextValue $ws.Range("D19") "0.3505"
Set-TextValue $ws.Range("E19") "0.57%"
Set-TextValue $ws.Range("D20") "8.761"
Set-TextValue $ws.Range("E20") "6.46%"
Set-TextValue $ws.Range("D21") "0.1363"
Set-TextValue $ws.Range("E21") "-0.14%"
Set-TextValue $ws.Range("D23") "0.04432"
Set-TextValue $ws.Range("E23") "0.85%"
Set-TextValue $ws.Range("D24") "0.001262"
Set-TextValue $ws.Range("E24") "0.49%"
Set-TextValue $ws.Range("E25") "0.94%"
Set-TextValue $ws.Range("D26") "0.0001263"
Set-TextValue $ws.Range("E26") "5.11%"
Set-TextValue $ws.Range("D27") "0.0003996"
Set-TextValue $ws.Range("D39") "0.02806"
Set-TextValue $ws.Range("E39") "1.12%"
Set-TextValue $ws.Range("D40") "0.05799"
Set-TextValue $ws.Range("E40") "7.50%"
Set-TextValue $ws.Range("D41") "0.007946"
Set-TextValue $ws.Range("E41") "2.17%"
Set-TextValue $ws.Range("D42") "0.1429"
Set-TextValue $ws.Range("E42") "0.89%"
Set-TextValue $ws.Range("D43") "0.008979"
Set-TextValue $ws.Range("E43") "0.34%"
Set-TextValue $ws.Range("E44") "0.05%"
Set-TextValue $ws.Range("D45") "0.01044"
Set-TextValue $ws.Range("E45") "-9.41%"
Set-TextValue $ws.Range("D46") "0.00007293"
Set-TextValue $ws.Range("E46") "7.84%"
Set-TextValue $ws.Range("E47") "0.11%"
Set-TextValue $ws.Range("D48") "0.003191"
Set-TextValue $ws.Range("E48") "-0.32%"
Set-TextValue $ws.Range("D49") "0.002274"
Set-TextValue $ws.Range("E49") "-0.32%"
Set-TextValue $ws.Range("D50") "0.00002104"
Set-TextValue $ws.Range("E50") "0.11%"
Set-TextValue $ws.Range("D51") "0.0002004"
Set-TextValue $ws.Range("E51") "0.11%"
